# Apply the update to the repayment report workbook.
#
# Summary of the change:
#   - Worksheet renamed from "repayment_20250901_20250915 (2)"
#     to "repayment_20250901_20250915 (4)"
#   - Three collectors got revised Repayment_collections (D),
#     Repayment_amount (E) and Pending Amount Recovery (G) figures:
#       Row 5  (Yandi Nugraha)        D 42 -> 43, E 24,063,135.00 -> 24,163,496.00, G 7.84 -> 7.87
#       Row 8  (Annisa Putri Restu)   D 40 -> 41, E 26,506,306.00 -> 27,192,618.00, G 7.92 -> 8.13
#       Row 14 (Nur Halim)            D 38 -> 43, E 25,023,925.00 -> 34,482,717.00, G 7.72 -> 10.64

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # The amount/rate columns hold numeric-looking values that are stored as
    # text in the workbook (e.g. "24,163,496.00"). Assigning a plain string to
    # .Value lets Excel auto-convert it to a number, so we momentarily force a
    # text number format while writing the value, then restore the cell style
    # so no extra formatting is left behind on the cell.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 5: Yandi Nugraha
$ws.Range("D5").Value = 43
Set-TextValue "E5" "24,163,496.00"
Set-TextValue "G5" "7.87"

# Row 8: Annisa Putri Restu
$ws.Range("D8").Value = 41
Set-TextValue "E8" "27,192,618.00"
Set-TextValue "G8" "8.13"

# Row 14: Nur Halim
$ws.Range("D14").Value = 43
Set-TextValue "E14" "34,482,717.00"
Set-TextValue "G14" "10.64"

# Rename the sheet to reflect the new upload revision.
$ws.Name = "repayment_20250901_20250915 (4)"

"Applied updates to rows 5, 8 and 14; renamed sheet to '$($ws.Name)'."
